$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update description text in E7 and E8 (swap order happened in shared strings but text content changes)
$ws.Range("E8").Value = "Creating + implementing splash screen, creating higher resolution icons + bordering them, class titles with proper font"
$ws.Range("E7").Value = "Routing, styling, creating and invoking existing classes based on each World of Warcraft class"

# Update hours worked on row 8 (end time 15 -> 16)
$ws.Range("C8").Value = 16

# Fill in row 9 (new work session entry)
$ws.Range("B9").Value = 10
$ws.Range("C9").Value = 16
$ws.Range("E9").Value = "Looking into noSQL options for this app and fending off error after error trying to implement them. Did not get much actual work done. "

# Update the selected cell/active selection on the sheet
$ws.Range("E9").Select()

$wb.Save()
